$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C contains the "Förändrad" (Changed) date, which is bumped by one
# day (46060 -> 46061) for every data row (rows 2 through 518).
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46060) {
        $cell.Value2 = 46061
    }
}
